$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.5
$ws.Range("H2").Value = 2.88
$ws.Range("K2").Value = 1.8
$ws.Range("Q2").Value = 2.33
$ws.Range("R2").Value = 1.61
$ws.Range("T2").Value = 1.36
$ws.Range("U2").Value = 5.6

# Row 3
$ws.Range("G3").Value = 2.2
$ws.Range("H3").Value = 2.9
$ws.Range("K3").Value = 1.91
$ws.Range("T3").Value = 1.44

# Row 4
$ws.Range("K4").Value = 1.83

# Row 5
$ws.Range("K5").Value = 1.83

# Row 6
$ws.Range("J6").Value = 1.95

# Row 7
$ws.Range("G7").Value = 1.72
$ws.Range("M7").Value = 1.03
$ws.Range("O7").Value = 1.19
$ws.Range("P7").Value = 4.33
$ws.Range("W7").Value = 2.7
$ws.Range("X7").Value = 1.41
$ws.Range("AC7").Value = 8.5
$ws.Range("AQ7").Value = 29
$ws.Range("AS7").Value = 151

# Row 8
$ws.Range("G8").Value = 2.2
$ws.Range("J8").Value = 2.88
$ws.Range("L8").Value = 3.6
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("S8").Value = 1.93
$ws.Range("T8").Value = 1.93
$ws.Range("W8").Value = 3.25
$ws.Range("X8").Value = 1.33
$ws.Range("Y8").Value = 1.4
$ws.Range("Z8").Value = 2.75
$ws.Range("AE8").Value = 9.5
$ws.Range("AI8").Value = 10
$ws.Range("AQ8").Value = 23
$ws.Range("AS8").Value = 201

# Row 9
$ws.Range("G9").Value = 1.55
$ws.Range("K9").Value = 2.38

# Row 13
$ws.Range("I13").Value = 2.1

# Row 14
$ws.Range("G14").Value = 3.8
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 1.91
$ws.Range("J14").Value = 4.33
$ws.Range("L14").Value = 2.5
$ws.Range("S14").Value = 1.85
$ws.Range("T14").Value = 2
$ws.Range("Y14").Value = 1.36
$ws.Range("Z14").Value = 3
$ws.Range("AD14").Value = 19
$ws.Range("AE14").Value = 13
$ws.Range("AG14").Value = 29
$ws.Range("AJ14").Value = 7
$ws.Range("AM14").Value = 7.5
$ws.Range("AN14").Value = 9.5
$ws.Range("AO14").Value = 8.5
$ws.Range("AP14").Value = 17
$ws.Range("AQ14").Value = 15

# Row 15
$ws.Range("I15").Value = 1.73
$ws.Range("S15").Value = 2.3
$ws.Range("T15").Value = 1.6
$ws.Range("W15").Value = 4.33
$ws.Range("X15").Value = 1.2

# Row 16
$ws.Range("Q16").Value = 1.83
$ws.Range("R16").Value = 1.98

# Row 17
$ws.Range("G17").Value = 2.47
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 2.52
$ws.Range("J17").Value = 2.85
$ws.Range("K17").Value = 2.37
$ws.Range("L17").Value = 2.95
$ws.Range("Z17").Value = 3.65
$ws.Range("AA17").Value = 1.35
$ws.Range("AB17").Value = 2.92
$ws.Range("AC17").Value = 16
$ws.Range("AD17").Value = 19.5
$ws.Range("AG17").Value = 17
$ws.Range("AH17").Value = 17
$ws.Range("AJ17").Value = 8
$ws.Range("AK17").Value = 10
$ws.Range("AL17").Value = 26
$ws.Range("AM17").Value = 14.5
$ws.Range("AN17").Value = 17.5
$ws.Range("AQ17").Value = 17.5
$ws.Range("AR17").Value = 18
$ws.Range("AS17").Value = 110

# Row 18
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 1.6
$ws.Range("J18").Value = 4.8
$ws.Range("L18").Value = 2.1
$ws.Range("P18").Value = 4.35
$ws.Range("W18").Value = 2.32
$ws.Range("X18").Value = 1.55
$ws.Range("Z18").Value = 3.2
$ws.Range("AB18").Value = 2.2
$ws.Range("AC18").Value = 18.5
$ws.Range("AD18").Value = 35
$ws.Range("AE18").Value = 15.5
$ws.Range("AF18").Value = 90
$ws.Range("AG18").Value = 40
$ws.Range("AJ18").Value = 8.25
$ws.Range("AK18").Value = 13.5
$ws.Range("AN18").Value = 8.75
$ws.Range("AP18").Value = 12.5
$ws.Range("AQ18").Value = 11.5
$ws.Range("AS18").Value = 300

# Row 19
$ws.Range("G19").Value = 2.5
$ws.Range("H19").Value = 2.88
$ws.Range("I19").Value = 3.2
$ws.Range("K19").Value = 1.91
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 7
$ws.Range("Q19").Value = 1.9
$ws.Range("R19").Value = 1.95
$ws.Range("AC19").Value = 6.5
$ws.Range("AF19").Value = 23
$ws.Range("AG19").Value = 23
$ws.Range("AM19").Value = 8
$ws.Range("AN19").Value = 15
$ws.Range("AP19").Value = 34

# Row 20
$ws.Range("H20").Value = 4.35
$ws.Range("I20").Value = 1.5
$ws.Range("J20").Value = 5.3
$ws.Range("K20").Value = 2.4
$ws.Range("L20").Value = 1.98
$ws.Range("M20").Value = 1.03
$ws.Range("N20").Value = 9
$ws.Range("O20").Value = 1.19
$ws.Range("P20").Value = 4.15
$ws.Range("S20").Value = 1.57
$ws.Range("T20").Value = 2.25
$ws.Range("W20").Value = 2.4
$ws.Range("X20").Value = 1.5
$ws.Range("Y20").Value = 1.31
$ws.Range("Z20").Value = 3.2
$ws.Range("AA20").Value = 1.72
$ws.Range("AB20").Value = 2.02
$ws.Range("AC20").Value = 18.5
$ws.Range("AD20").Value = 37
$ws.Range("AG20").Value = 50
$ws.Range("AH20").Value = 45
$ws.Range("AI20").Value = 9
$ws.Range("AJ20").Value = 8.5
$ws.Range("AL20").Value = 60
$ws.Range("AM20").Value = 8.25
$ws.Range("AO20").Value = 8.25
$ws.Range("AP20").Value = 10.75
$ws.Range("AS20").Value = 400
